$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4101.2
$ws.Range("I86").Value = 4101.2
$ws.Range("K86").Value = 4101.2
$ws.Range("M86").Value = -2978.2

$ws.Range("H89").Value = 4101.2
$ws.Range("I89").Value = 4101.2
$ws.Range("K89").Value = 20506
$ws.Range("M89").Value = -14890

$ws.Range("H100").Value = 2182.6667
$ws.Range("I100").Value = 2182.6667
$ws.Range("K100").Value = 2182.6667
$ws.Range("M100").Value = -1641.6667

$ws.Range("H113").Value = 34089.5
$ws.Range("I113").Value = 33850.668
$ws.Range("J113").Value = 34328.332
$ws.Range("K113").Value = 33850.668
$ws.Range("L113").Value = 34328.332
$ws.Range("M113").Value = -30596.668
$ws.Range("N113").Value = -40836.332

$ws.Range("H132").Value = 4000.75
$ws.Range("J132").Value = 5499.5
$ws.Range("L132").Value = 16498.5
$ws.Range("N132").Value = -21558.5

$ws.Range("H137").Value = 5331.6665
$ws.Range("J137").Value = 2996
$ws.Range("L137").Value = 8988
$ws.Range("N137").Value = -14088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H41").Value = 1998
$ws.Range("I41").Value = 1998
$ws.Range("K41").Value = 1998
$ws.Range("M41").Value = -1584

$ws.Range("H92").Value = 57500
$ws.Range("J92").Value = 57500
$ws.Range("L92").Value = 57500
$ws.Range("N92").Value = -62492

$ws.Range("H104").Value = 57056.25
$ws.Range("J104").Value = 57056.25
$ws.Range("L104").Value = 57056.25
$ws.Range("N104").Value = -64044.25

$ws.Range("H132").Value = 4499.6665
$ws.Range("I132").Value = 4499.6665
$ws.Range("K132").Value = 13498.9995
$ws.Range("M132").Value = -10968.9995

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5454

$ws.Range("H86").Value = 1499
$ws.Range("I86").Value = 1499
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1499
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -376
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 1499
$ws.Range("I89").Value = 1499
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 7495
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1879
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 2391.7273
$ws.Range("I94").Value = 1968.1666
$ws.Range("K94").Value = 1968.1666
$ws.Range("M94").Value = -1517.1666

$ws.Range("H105").Value = 2010
$ws.Range("I105").Value = 2010
$ws.Range("K105").Value = 2010
$ws.Range("M105").Value = -263

$ws.Range("H107").Value = 1408.4286
$ws.Range("I107").Value = 1408.4286
$ws.Range("K107").Value = 1408.4286
$ws.Range("M107").Value = 511.5714

$ws.Range("H134").Value = 5999.619
$ws.Range("I134").Value = 5299.6
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 15898.8
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -13363.8
$ws.Range("N134").Value = -65070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1076.6364
$ws.Range("I12").Value = 675
$ws.Range("K12").Value = 675
$ws.Range("M12").Value = -505

$ws.Range("H22").Value = 801
$ws.Range("I22").Value = 801
$ws.Range("K22").Value = 801
$ws.Range("M22").Value = -451

$ws.Range("H31").Value = 9811.875
$ws.Range("I31").Value = 9099.200000000001
$ws.Range("K31").Value = 9099.200000000001
$ws.Range("M31").Value = -8804.200000000001

$ws.Range("H34").Value = 9811.875
$ws.Range("I34").Value = 9099.200000000001
$ws.Range("K34").Value = 9099.200000000001
$ws.Range("M34").Value = -8897.200000000001

$ws.Range("H35").Value = 5380
$ws.Range("J35").Value = 5274
$ws.Range("L35").Value = 5274
$ws.Range("N35").Value = -5862

$ws.Range("H58").Value = 992.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 992.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 992.5
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -1398.5

$ws.Range("H132").Value = 2508.1667
$ws.Range("I132").Value = 2233.75
$ws.Range("K132").Value = 6701.25
$ws.Range("M132").Value = -4171.25

$ws.Range("H136").Value = 992.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 992.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 2977.5
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -8077.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 230.66667
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 345
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 1035
$ws.Range("M12").Value = 167
$ws.Range("N12").Value = -1381

$ws.Range("H18").Value = 694.5
$ws.Range("I18").Value = 819.5
$ws.Range("J18").Value = 444.5
$ws.Range("K18").Value = 2458.5
$ws.Range("L18").Value = 1333.5
$ws.Range("M18").Value = -2289.5
$ws.Range("N18").Value = -1671.5

$ws.Range("H129").Value = 3387.5
$ws.Range("I129").Value = 375
$ws.Range("K129").Value = 1125
$ws.Range("M129").Value = 3875

$ws.Range("H131").Value = 2306.2
$ws.Range("J131").Value = 3010.8333
$ws.Range("L131").Value = 9032.499899999999
$ws.Range("N131").Value = -19112.4999

$ws.Range("H136").Value = 1030
$ws.Range("I136").Value = 1030
$ws.Range("K136").Value = 3090
$ws.Range("M136").Value = 2010

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5212.5
$ws.Range("I80").Value = 4425
$ws.Range("K80").Value = 4425
$ws.Range("M80").Value = -3427

$ws.Range("H83").Value = 5212.5
$ws.Range("I83").Value = 4425
$ws.Range("K83").Value = 22125
$ws.Range("M83").Value = -17133

$ws.Range("H104").Value = 209835.5
$ws.Range("J104").Value = 209835.5
$ws.Range("L104").Value = 209835.5
$ws.Range("N104").Value = -216823.5

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H122").Value = 5002.3335
$ws.Range("I122").Value = 5002.3335
$ws.Range("K122").Value = 15007.0005
$ws.Range("M122").Value = -12557.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1201
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -888

$ws.Range("H17").Value = 8099.857
$ws.Range("I17").Value = 3349.5
$ws.Range("J17").Value = 10000
$ws.Range("K17").Value = 3349.5
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = -3179.5
$ws.Range("N17").Value = -10340

$ws.Range("H39").Value = 60000
$ws.Range("J39").Value = 60000
$ws.Range("L39").Value = 60000
$ws.Range("N39").Value = -60920

$ws.Range("H40").Value = 52002.4
$ws.Range("I40").Value = 40002.332
$ws.Range("J40").Value = 70002.5
$ws.Range("K40").Value = 40002.332
$ws.Range("L40").Value = 70002.5
$ws.Range("M40").Value = -39866.332
$ws.Range("N40").Value = -70274.5

$ws.Range("H61").Value = 100005
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 100005
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 100005
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -100409

$ws.Range("H113").Value = 100005
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 100005
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 100005
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -104345

$ws.Range("H126").Value = 1201
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8664.666999999999
$ws.Range("I2").Value = 19995
$ws.Range("K2").Value = 19995
$ws.Range("M2").Value = -19883

$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939

$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696

$ws.Range("H113").Value = 7463
$ws.Range("J113").Value = 18649.334
$ws.Range("L113").Value = 55948.00199999999
$ws.Range("N113").Value = -60288.00199999999

$ws.Range("H132").Value = 4450
$ws.Range("I132").Value = 4450
$ws.Range("K132").Value = 13350
$ws.Range("M132").Value = -10820

$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

$ws.Range("H140").Value = 52500
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360
